# Penalty Reward System (unfinished) - shift forecast window forward by one
# week on the "Forecast Comparison" sheet and refresh a handful of derived
# figures on the "Summary" sheet.
#
# Note: date-shaped and pure-numeric strings get auto-coerced to real
# dates/numbers by the Value setter, so they are written with a leading
# apostrophe to force plain text, matching the workbook's original
# (inlineStr) cell typing.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$weekStartDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$myForecasts = @(19, 20, 21, 21, 16, 16, 16, 16, 16, 16, 16, 20, 16, 16, 16, 16)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $wsForecast.Range("B$row").Value = "'" + $weekStartDates[$i]
    $wsForecast.Range("D$row").Value = $myForecasts[$i]
}

$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = "2024-02-11 to 2025-01-05"
$wsSummary.Range("B5").Value = "'28"
$wsSummary.Range("B9").Value = "'278"
$wsSummary.Range("B10").Value = "'144"
$wsSummary.Range("B11").Value = "'80"
$wsSummary.Range("B12").Value = "'21"
$wsSummary.Range("B13").Value = "'2025-01-26"
$wsSummary.Range("B14").Value = "'16"
$wsSummary.Range("B15").Value = "'2025-02-09"
